$d = $word.ActiveDocument

# 1) Update w:lang on the picture-bearing runs: ru-RU -> uk-UA (both w:val and w:eastAsia)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $rr = $p.Range
    if ($rr.LanguageIDFarEast -eq "ru-RU") {
        $rr.LanguageID = "uk-UA"
        $rr.LanguageIDFarEast = "uk-UA"
    }
}

# 2) Drop the paragraph-mark rPr/lang override on the "7. Запустите AutoCAD..." paragraph
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("7. ")) {
        $p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00CE7ECA" w:rsidRDefault="00CE7ECA"><w:r w:rsidRPr="00CE7ECA"><w:lastRenderedPageBreak/><w:t xml:space="preserve">7. </w:t></w:r><w:r><w:t xml:space="preserve">Запустите </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>AutoCAD</w:t></w:r><w:r w:rsidRPr="00CE7ECA"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>на выполнение</w:t></w:r></w:p>')
        break
    }
}

# 3) Fix the typo / move the _GoBack bookmark / drop proofErr markers in item 9
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.Contains("Используйте команду")) {
        $p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00DE297E" w:rsidRPr="00DE297E" w:rsidRDefault="00DE297E"><w:r><w:t xml:space="preserve">9. Используйте команду </w:t></w:r><w:r w:rsidRPr="00DE297E"><w:rPr><w:lang w:val="el-GR"/></w:rPr><w:t>R</w:t></w:r><w:r w:rsidRPr="00DE297E"><w:t>-</w:t></w:r><w:r w:rsidRPr="00DE297E"><w:rPr><w:lang w:val="el-GR"/></w:rPr><w:t>TRIANG</w:t></w:r><w:r><w:t xml:space="preserve"> для построения развертки ли</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>нейчатой поверхности.</w:t></w:r></w:p>')
        break
    }
}
